# Commit: "add B5, B6, and B7"
# Fill in the real TxHash evidence values for sheets B5, B6 and B7 (A2/A3),
# which previously all shared the same generic placeholder labels, and
# extend each sheet's used range down to row 6 (blank rows 5 & 6), matching
# the sibling sheets (e.g. B1) that already carry extra trailing rows.
#
# Order matters for shared-string slot packing: write B5 first, then B6,
# then B7 so the now-unreferenced placeholder strings get reclaimed and the
# new unique strings land in the same slot order as the target workbook.

$wb = $excel.ActiveWorkbook

$ws5 = $wb.Worksheets.Item("B5")
$ws5.Range("A2").Value = "A28D4C47D418CCCE69F8CFD1F91B9CB3FDD5AA40D12B390648A0217F1E1C1668"
$ws5.Range("A3").Value = "9C5AC0E16C034D1519D2139A7A3A32CD10C7D3E4C4DC17570E6A72A3B2E2C10F"
$ws5.Range("A5").Value = "x"
$ws5.Range("A6").Value = "x"
$ws5.Range("A5:A6").Value = ""

$ws6 = $wb.Worksheets.Item("B6")
$ws6.Range("A2").Value = "B8D72D04D01CFE297F3E1ED194FA9E40EA6E2FA6D0E757F44E52484F34E29117"
$ws6.Range("A3").Value = "1749A47D641A9665DA61FB1C0C877E87046FE4636015264235CD403EE63C09B7"
$ws6.Range("A5").Value = "x"
$ws6.Range("A6").Value = "x"
$ws6.Range("A5:A6").Value = ""

$ws7 = $wb.Worksheets.Item("B7")
$ws7.Range("A2").Value = "A57724760F7669D776EADE6F39164F3F699F8A8A1D5D9B77D626F87BE27F4818"
$ws7.Range("A3").Value = "33E160228043C0C6658E9BD18CCD0AE89694AD4B160D4F9EC357561953B9F063"
$ws7.Range("A5").Value = "x"
$ws7.Range("A6").Value = "x"
$ws7.Range("A5:A6").Value = ""

# Restore per-sheet selections, then land on B5 as the active tab/cell
# (matches the final saved cursor position in the workbook).
$ws6.Range("B6").Select()
$ws7.Range("A5").Select()

$ws5.Activate()
$ws5.Range("I22").Select()

Write-Output "done"
